$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 19.14999961853028
$ws.Range("E2").Value = 17.98999977111816
$ws.Range("F2").Value = 19.23999977111816
$ws.Range("G2").Value = 17.68000030517578
$ws.Range("H2").Value = 272490000
$ws.Range("I2").Value = "CDNS"

$ws.Range("D3").Value = 18.44000053405762
$ws.Range("E3").Value = 18.64999961853028
$ws.Range("F3").Value = 19.82999992370605
$ws.Range("G3").Value = 17.77000045776367
$ws.Range("H3").Value = 272490000
$ws.Range("I3").Value = "CDNS"

$ws.Range("D4").Value = 19.80999946594238
$ws.Range("E4").Value = 20.96999931335449
$ws.Range("F4").Value = 21.13999938964844
$ws.Range("G4").Value = 18.8799991607666
$ws.Range("H4").Value = 272490000
$ws.Range("I4").Value = "CDNS"

$ws.Range("D5").Value = 20.73999977111816
$ws.Range("E5").Value = 22.21999931335449
$ws.Range("F5").Value = 23.29999923706055
$ws.Range("G5").Value = 20.07999992370605
$ws.Range("H5").Value = 272490000
$ws.Range("I5").Value = "CDNS"

$ws.Range("D6").Value = 20
$ws.Range("E6").Value = 19.55999946594238
$ws.Range("F6").Value = 20.63999938964844
$ws.Range("G6").Value = 18.31999969482422
$ws.Range("H6").Value = 272490000
$ws.Range("I6").Value = "CDNS"

$ws.Range("D7").Value = 23.48999977111816
$ws.Range("E7").Value = 23.19000053405762
$ws.Range("F7").Value = 23.88999938964844
$ws.Range("G7").Value = 22.79000091552734
$ws.Range("H7").Value = 272490000
$ws.Range("I7").Value = "CDNS"

$ws.Range("D8").Value = 24.29999923706055
$ws.Range("E8").Value = 24.04999923706055
$ws.Range("F8").Value = 26.23999977111816
$ws.Range("G8").Value = 23.82999992370605
$ws.Range("H8").Value = 272490000
$ws.Range("I8").Value = "CDNS"

$ws.Range("D9").Value = 25.5
$ws.Range("E9").Value = 25.57999992370605
$ws.Range("F9").Value = 28
$ws.Range("G9").Value = 24.98999977111816
$ws.Range("H9").Value = 272490000
$ws.Range("I9").Value = "CDNS"

$ws.Range("D10").Value = 25.43000030517578
$ws.Range("E10").Value = 26.03000068664551
$ws.Range("F10").Value = 26.46999931335449
$ws.Range("G10").Value = 25.23999977111816
$ws.Range("H10").Value = 272490000
$ws.Range("I10").Value = "CDNS"

$ws.Range("D11").Value = 31.39999961853028
$ws.Range("E11").Value = 32.56999969482422
$ws.Range("F11").Value = 33.61000061035156
$ws.Range("G11").Value = 30.80999946594238
$ws.Range("H11").Value = 272490000
$ws.Range("I11").Value = "CDNS"

$ws.Range("D12").Value = 33.58000183105469
$ws.Range("E12").Value = 36.90000152587891
$ws.Range("F12").Value = 37.5099983215332
$ws.Range("G12").Value = 32.84999847412109
$ws.Range("H12").Value = 272490000
$ws.Range("I12").Value = "CDNS"

$ws.Range("D13").Value = 39.79000091552734
$ws.Range("E13").Value = 43.15999984741211
$ws.Range("F13").Value = 43.22999954223633
$ws.Range("G13").Value = 39.59999847412109
$ws.Range("H13").Value = 272490000
$ws.Range("I13").Value = "CDNS"

$ws.Range("D14").Value = 42
$ws.Range("E14").Value = 44.86000061035156
$ws.Range("F14").Value = 46
$ws.Range("G14").Value = 41.70999908447266
$ws.Range("H14").Value = 272490000
$ws.Range("I14").Value = "CDNS"

$ws.Range("D15").Value = 36.68999862670898
$ws.Range("E15").Value = 40.06000137329102
$ws.Range("F15").Value = 41.4900016784668
$ws.Range("G15").Value = 35.4900016784668
$ws.Range("H15").Value = 272490000
$ws.Range("I15").Value = "CDNS"

$ws.Range("D16").Value = 42.95000076293945
$ws.Range("E16").Value = 44.09000015258789
$ws.Range("F16").Value = 46.9900016784668
$ws.Range("G16").Value = 42.84000015258789
$ws.Range("H16").Value = 272490000
$ws.Range("I16").Value = "CDNS"

$ws.Range("D17").Value = 45.45999908447266
$ws.Range("E17").Value = 44.56999969482422
$ws.Range("F17").Value = 47.04000091552734
$ws.Range("G17").Value = 39.08000183105469
$ws.Range("H17").Value = 272490000
$ws.Range("I17").Value = "CDNS"

$ws.Range("D18").Value = 42.65000152587891
$ws.Range("E18").Value = 48.02999877929688
$ws.Range("F18").Value = 48.16999816894531
$ws.Range("G18").Value = 41.43000030517578
$ws.Range("H18").Value = 272490000
$ws.Range("I18").Value = "CDNS"

$ws.Range("D19").Value = 65.08000183105469
$ws.Range("E19").Value = 69.37999725341797
$ws.Range("F19").Value = 69.69000244140625
$ws.Range("G19").Value = 62.52000045776367
$ws.Range("H19").Value = 272490000
$ws.Range("I19").Value = "CDNS"

$ws.Range("D20").Value = 72.55999755859375
$ws.Range("E20").Value = 73.91000366210938
$ws.Range("F20").Value = 77.08000183105469
$ws.Range("G20").Value = 71.09999847412109
$ws.Range("H20").Value = 272490000
$ws.Range("I20").Value = "CDNS"

$ws.Range("D21").Value = 66.34999847412109
$ws.Range("E21").Value = 65.34999847412109
$ws.Range("F21").Value = 70.5
$ws.Range("G21").Value = 62.40999984741211
$ws.Range("H21").Value = 272490000
$ws.Range("I21").Value = "CDNS"

$ws.Range("D22").Value = 70.18000030517578
$ws.Range("E22").Value = 72.11000061035156
$ws.Range("F22").Value = 76.30999755859375
$ws.Range("G22").Value = 69.58999633789062
$ws.Range("H22").Value = 272490000
$ws.Range("I22").Value = "CDNS"

$ws.Range("D23").Value = 63.43000030517578
$ws.Range("E23").Value = 81.12999725341797
$ws.Range("F23").Value = 82.77999877929688
$ws.Range("G23").Value = 62.91999816894531
$ws.Range("H23").Value = 272490000
$ws.Range("I23").Value = "CDNS"

$ws.Range("D24").Value = 95.6999969482422
$ws.Range("E24").Value = 109.25
$ws.Range("F24").Value = 109.2600021362305
$ws.Range("G24").Value = 95.23000335693359
$ws.Range("H24").Value = 272490000
$ws.Range("I24").Value = "CDNS"

$ws.Range("D25").Value = 107.8499984741211
$ws.Range("E25").Value = 109.370002746582
$ws.Range("F25").Value = 118.2200012207031
$ws.Range("G25").Value = 104.4800033569336
$ws.Range("H25").Value = 272490000
$ws.Range("I25").Value = "CDNS"

$ws.Range("D26").Value = 136.8899993896484
$ws.Range("E26").Value = 130.3899993896484
$ws.Range("F26").Value = 142.1100006103516
$ws.Range("G26").Value = 127.5599975585938
$ws.Range("H26").Value = 272490000
$ws.Range("I26").Value = "CDNS"

$ws.Range("D27").Value = 138.4600067138672
$ws.Range("E27").Value = 131.7700042724609
$ws.Range("F27").Value = 148.2700042724609
$ws.Range("G27").Value = 130.2700042724609
$ws.Range("H27").Value = 272490000
$ws.Range("I27").Value = "CDNS"

$ws.Range("D28").Value = 136.0399932861328
$ws.Range("E28").Value = 147.6499938964844
$ws.Range("F28").Value = 148.1199951171875
$ws.Range("G28").Value = 135.0399932861328
$ws.Range("H28").Value = 272490000
$ws.Range("I28").Value = "CDNS"

$ws.Range("D29").Value = 152.3399963378906
$ws.Range("E29").Value = 173.1100006103516
$ws.Range("F29").Value = 175.0800018310547
$ws.Range("G29").Value = 145.8399963378906
$ws.Range("H29").Value = 272490000
$ws.Range("I29").Value = "CDNS"

$ws.Range("D30").Value = 185.6799926757812
$ws.Range("E30").Value = 152.1399993896484
$ws.Range("F30").Value = 188.5700073242188
$ws.Range("G30").Value = 136.6199951171875
$ws.Range("H30").Value = 272490000
$ws.Range("I30").Value = "CDNS"

$ws.Range("D31").Value = 165
$ws.Range("E31").Value = 150.8500061035156
$ws.Range("F31").Value = 168.7299957275391
$ws.Range("G31").Value = 147.2100067138672
$ws.Range("H31").Value = 272490000
$ws.Range("I31").Value = "CDNS"

$ws.Range("D32").Value = 149.7599945068359
$ws.Range("E32").Value = 186.0800018310547
$ws.Range("F32").Value = 186.5800018310547
$ws.Range("G32").Value = 147
$ws.Range("H32").Value = 272490000
$ws.Range("I32").Value = "CDNS"

$ws.Range("D33").Value = 164.3899993896484
$ws.Range("E33").Value = 151.3899993896484
$ws.Range("F33").Value = 174.3300018310547
$ws.Range("G33").Value = 142.3399963378906
$ws.Range("H33").Value = 272490000
$ws.Range("I33").Value = "CDNS"

$ws.Range("D34").Value = 162.8300018310547
$ws.Range("E34").Value = 182.8300018310547
$ws.Range("F34").Value = 186.2400054931641
$ws.Range("G34").Value = 154.8800048828125
$ws.Range("H34").Value = 272490000
$ws.Range("I34").Value = "CDNS"

$ws.Range("D35").Value = 209
$ws.Range("E35").Value = 209.4499969482422
$ws.Range("F35").Value = 217.6699981689453
$ws.Range("G35").Value = 194.009994506836
$ws.Range("H35").Value = 272490000
$ws.Range("I35").Value = "CDNS"

$ws.Range("D36").Value = 234.3800048828125
$ws.Range("E36").Value = 234.009994506836
$ws.Range("F36").Value = 248.1600036621093
$ws.Range("G36").Value = 226.1100006103516
$ws.Range("H36").Value = 272490000
$ws.Range("I36").Value = "CDNS"

$ws.Range("D37").Value = 235.4400024414062
$ws.Range("E37").Value = 239.8500061035156
$ws.Range("F37").Value = 255.8500061035156
$ws.Range("G37").Value = 227.7700042724609
$ws.Range("H37").Value = 272490000
$ws.Range("I37").Value = "CDNS"

$ws.Range("D38").Value = 269.1600036621094
$ws.Range("E38").Value = 288.4599914550781
$ws.Range("F38").Value = 301.6099853515625
$ws.Range("G38").Value = 251.9400024414062
$ws.Range("H38").Value = 272490000
$ws.Range("I38").Value = "CDNS"

$ws.Range("D39").Value = 310.4599914550781
$ws.Range("E39").Value = 275.6300048828125
$ws.Range("F39").Value = 317.0899963378906
$ws.Range("G39").Value = 273.489990234375
$ws.Range("H39").Value = 272490000
$ws.Range("I39").Value = "CDNS"

$ws.Range("D40").Value = 307.75
$ws.Range("E40").Value = 267.6600036621094
$ws.Range("F40").Value = 322.2799987792969
$ws.Range("G40").Value = 250.8399963378907
$ws.Range("H40").Value = 272490000
$ws.Range("I40").Value = "CDNS"

$ws.Range("D41").Value = 273.4400024414062
$ws.Range("E41").Value = 276.1199951171875
$ws.Range("F41").Value = 290.2000122070312
$ws.Range("G41").Value = 246.6000061035156
$ws.Range("H41").Value = 272490000
$ws.Range("I41").Value = "CDNS"

$ws.Range("D42").Value = 304.8399963378906
$ws.Range("E42").Value = 297.6199951171875
$ws.Range("F42").Value = 325.0299987792969
$ws.Range("G42").Value = 285.0599975585937
$ws.Range("H42").Value = 272490000
$ws.Range("I42").Value = "CDNS"

$ws.Range("D43").Value = 252.6399993896484
$ws.Range("E43").Value = 297.739990234375
$ws.Range("F43").Value = 303.6000061035156
$ws.Range("G43").Value = 221.5599975585937
$ws.Range("H43").Value = 272490000
$ws.Range("I43").Value = "CDNS"

$ws.Range("D44").Value = 306.3500061035156
$ws.Range("E44").Value = 364.5700073242188
$ws.Range("F44").Value = 376.4500122070313
$ws.Range("G44").Value = 304.7300109863281
$ws.Range("H44").Value = 272490000
$ws.Range("I44").Value = "CDNS"
